# Applies the "Finished adding all the major parts" edit:
#  - inserts 4 rows at the top of Sheet1 for a new "MAJOR PARTS" banner
#  - shifts the existing drawings/hyperlinks to match the new row numbers
#  - appends a new "Binding Posts:" section (3 new links) after "USB Connectors"
#  - tweaks the sheet view (zoom / selection)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1. Insert 4 blank rows above everything else; all existing cell content,
#    styles and shared strings shift down by 4 rows automatically.
# ---------------------------------------------------------------------------
$ws.Rows("1:4").Insert()

# ---------------------------------------------------------------------------
# 2. Re-anchor the 3 pictures: Rows.Insert() does not move floating drawings,
#    so nudge each one down by the height of the 4 new rows (58pt @ the
#    sheet's default 14.5pt row height).
# ---------------------------------------------------------------------------
for ($i = 1; $i -le $ws.Shapes.Count; $i++) {
    $shp = $ws.Shapes.Item($i)
    $shp.Top = $shp.Top + 58
}

# ---------------------------------------------------------------------------
# 3. New title banner in row 1 (merges visually across A:B via matching
#    formatting), rows 2-3 are blank spacer rows with the same formatting.
# ---------------------------------------------------------------------------
$ws.Range("A1").Value = "MAJOR PARTS"
$titleRange = $ws.Range("A1:B3")
$titleRange.Font.Bold = $true
$titleRange.Font.Size = 26
$titleRange.VerticalAlignment = -4108   # xlVAlignCenter

$ws.Rows(1).RowHeight = 40.5
$ws.Rows(2).RowHeight = 14.5
$ws.Rows(3).RowHeight = 14.5

# ---------------------------------------------------------------------------
# 4. Append the new "Binding Posts:" block under "USB Connectors" (now on
#    row 48 after the insert above).
# ---------------------------------------------------------------------------
$ws.Range("A50").Value = "Binding Posts:"
$ws.Range("B50").Value = "https://www.digikey.ca/product-detail/en/sparkfun-electronics/PRT-09740/1568-1665-ND/7393681"
$ws.Range("B51").Value = "https://www.digikey.ca/product-detail/en/sparkfun-electronics/PRT-09739/1568-1664-ND/7393680"
$ws.Range("B53").Value = "https://www.ebay.ie/itm/DIY-PCB-Banana-jack-binding-post-breakout-board-/163248841920"

# ---------------------------------------------------------------------------
# 5. Rebuild the hyperlinks. Rows.Insert() left the old <hyperlinks> refs
#    pointing at the pre-shift addresses, so clear everything out and
#    re-create each link against its new cell, restoring the classic
#    "Hyperlink" cell style afterwards (Hyperlinks.Add otherwise stamps a
#    brand new style record even when one already fits).
# ---------------------------------------------------------------------------
$ws.Range("A1").Hyperlinks.Delete()

function Add-Link($addr, $url) {
    $ws.Hyperlinks.Add($ws.Range($addr), $url) | Out-Null
    $ws.Range($addr).Style = "Hyperlink"
}

Add-Link "A8"  "https://www.digikey.ca/product-detail/en/mean-well-usa-inc/IRM-60-24/1866-3067-ND/7704690"
Add-Link "A13" "https://www.digikey.ca/product-detail/en/linear-technology-analog-devices/LTC3649EFE-PBF/LTC3649EFE-PBF-ND/5825359"
Add-Link "A15" "https://www.digikey.ca/product-detail/en/linear-technology-analog-devices/LT3086EFE-TRPBF/LT3086EFE-TRPBFTR-ND/5233238"
Add-Link "B30" "https://www.digikey.ca/product-detail/en/bourns-inc/PTV09A-4020F-A103/PTV09A-4020F-A103-ND/3781119"
Add-Link "B31" "https://www.digikey.ca/product-detail/en/bourns-inc/PTV09A-4020F-A504/PTV09A-4020F-A504-ND/3781122"
Add-Link "A36" "https://www.digikey.ca/product-detail/en/recom-power/R-78E3.3-1.0/945-2409-5-ND/5327711"
Add-Link "A39" "https://www.digikey.ca/product-detail/en/recom-power/R-78B5.0-2.0/945-3042-ND/6677084"
Add-Link "B44" "http://www.electronics-diy.com/70v_pic_voltmeter_amperemeter.php"
Add-Link "B45" "https://www.adafruit.com/product/399?gclid=Cj0KCQjwjMfoBRDDARIsAMUjNZpkoKOfF8LGOWSxf4VLWIPIVsbuiiDFCeu1C8yUAKVEhcIb9canXt8aAnndEALw_wcB"
Add-Link "B46" "https://www.digikey.com/product-detail/en/PIC16F876A-I%2fSO/PIC16F876A-I%2fSO-ND/446139/?itemSeq=296474682"
Add-Link "B48" "https://www.digikey.ca/product-detail/en/te-connectivity-amp-connectors/1734366-1/A114947-ND/1891570"
Add-Link "B50" "https://www.digikey.ca/product-detail/en/sparkfun-electronics/PRT-09740/1568-1665-ND/7393681"
Add-Link "B51" "https://www.digikey.ca/product-detail/en/sparkfun-electronics/PRT-09739/1568-1664-ND/7393680"
Add-Link "B53" "https://www.ebay.ie/itm/DIY-PCB-Banana-jack-binding-post-breakout-board-/163248841920"

# ---------------------------------------------------------------------------
# 6. Sheet view tweaks: slightly smaller zoom, new active selection.
# ---------------------------------------------------------------------------
$ws.Activate()
$excel.ActiveWindow.Zoom = 85
$ws.Range("G5").Select()
